$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in columns D and E stay as text (not auto-converted to numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.485.95'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.660.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.27%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.23'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.95'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.32%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.614'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.128'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.07%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.397'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.83'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.154'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.16'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000194'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.143.44'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.373.98'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.48%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.680.72'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.76'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.76'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '349.99'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.56'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000110'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.55'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.62'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.57'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.97'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.98%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.12'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '520.26'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.75'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.39'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.41'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.52%  '

$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.49'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.63%  '

$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.420'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '157.60'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.92'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '162.41'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.09'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.28'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0603'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.57'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.637'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0256'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0260'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +13.53%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0995'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.92'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -5.56%  '
